$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Plain value edits on existing rows (no structural change) ---
# Row 7: num_controls corrected 71 -> 76 (ripples into J7/K7 automatically)
$ws.Range("C7").Value = 76

# Row 9: num_controls / num_cases corrected
$ws.Range("C9").Value = 31
$ws.Range("D9").Value = 10

# --- 2. Relocate the GSE57605 entry (currently row 12) further down the sheet ---
# Capture its current values before the row is removed from the main block.
$gse57605 = @{
    A = $ws.Range("A12").Value2
    B = $ws.Range("B12").Value2
    C = $ws.Range("C12").Value2
    D = $ws.Range("D12").Value2
    E = $ws.Range("E12").Value2
    F = $ws.Range("F12").Value2
    G = $ws.Range("G12").Value2
    H = $ws.Range("H12").Value2
    I = $ws.Range("I12").Value2
}

# Remove row 12 from the main table; rows 13:21 shift up to 12:20.
$ws.Rows("12").Delete()

# The row that used to be 13 (GSE63311) is now row 12 - refresh its edited figures.
$ws.Range("C12").Value = 46
$ws.Range("D12").Value = 37
$ws.Range("M12").Formula = "=83-37"

# M3 / M4 used to reference the old row 16 / row 21 totals - point them at the
# rows those totals now live in (15 / 20) after the shift.
$ws.Range("M3").Formula = "=K3+K2+K5+K6+K7+K8+K15"
$ws.Range("M4").Formula = "=M3/K20"

# The totals row's weighted-average formula (now row 20) dropped the deleted
# row's K12 term and needs its surviving references re-pointed post-shift.
$ws.Range("J20").Formula = "=SUM(J7:J18)/(K7+K10+K11+K14+K15+K16+K18)"

# --- 3. Re-create the relocated GSE57605 entry at row 24 ---
$ws.Range("A24").Value = $gse57605.A
$ws.Range("B24").Value = $gse57605.B
$ws.Range("C24").Value = $gse57605.C
$ws.Range("D24").Value = $gse57605.D
$ws.Range("E24").Value = $gse57605.E
$ws.Range("F24").Value = $gse57605.F
$ws.Range("G24").Value = $gse57605.G
$ws.Range("H24").Value = $gse57605.H
$ws.Range("I24").NumberFormat = "0%"
$ws.Range("I24").Value = $gse57605.I
$ws.Range("J24").Formula = "=I24*K24"
$ws.Range("K24").Formula = "=C24+D24"

# --- 4. New reconciliation row underneath it ---
$ws.Range("C28").Formula = "=463-25"
$ws.Range("D28").Formula = "=1232-D24"

# --- 5. Selection follows the new last data row ---
$ws.Range("A19").Select()
